$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header / note area (rows 1-3)
# ---------------------------------------------------------------
$ws.Range("A1").Value = "Note: character limit of about 45."

$ws.Range("D2").Value = "Work days vs. Saturday/Sunday detections."
$ws.Range("D2").Font.Bold = $true

$ws.Range("J2").Value = "Siri integration? What should I wear."
$ws.Range("J3").Value = "Target athletes?"

# ---------------------------------------------------------------
# "Extras" label near the table header row
# ---------------------------------------------------------------
$ws.Range("O8").Value = "Extras"
$ws.Range("O8").Font.Bold = $true

# ---------------------------------------------------------------
# Temp range label update: "0 - 32º" -> "32 - 49 "
# ---------------------------------------------------------------
$ws.Range("A11").Value = "32 - 49 "
$ws.Range("A12").Value = "32 - 49 "
$ws.Range("A13").Value = "32 - 49 "
$ws.Range("A14").Value = "32 - 49 "
$ws.Range("A15").Value = "32 - 49 "
$ws.Range("A16").Value = "32 - 49 "

# ---------------------------------------------------------------
# Clear out the stray B11:B16 cells (no longer used)
# ---------------------------------------------------------------
$ws.Range("B11").Clear()
$ws.Range("B12").Clear()
$ws.Range("B13").Clear()
$ws.Range("B14").Clear()
$ws.Range("B15").Clear()
$ws.Range("B16").Clear()

# ---------------------------------------------------------------
# New punchy phrases table, rows 11-16, columns C/D/E + merged F:G
# ---------------------------------------------------------------
$ws.Range("C11").Value = "A little rainy. Dress warm, boots ideal."
$ws.Range("D11").Value = "It's very cold and rainy. Heavy raincoat."
$ws.Range("E11").Value = "Super rainy. Just stay inside."
$ws.Range("F11").Value = "Umbrella would probably break. Big indoors day."

$ws.Range("C12").Value = "You won't get too wet. It's cold though."
$ws.Range("D12").Value = "Colder than a penguin's ass. Dress heavy."
$ws.Range("E12").Value = "Cold and rainy. Great day for some Netflix."
$ws.Range("F12").Value = "Stay inside if you're sane."

$ws.Range("C13").Value = "Pretty nippy. Layer up."
$ws.Range("D13").Value = "It's cold as hell and wet. Wear a coat."
$ws.Range("E13").Value = "Don't go out unless you love being drenched."
$ws.Range("F13").Value = "Get your dog inside and don't leave for a bit."

$ws.Range("C14").Value = "Not the worst, but it's kinda rainy."
$ws.Range("D14").Value = "Wear something water-resistant and warm."
$ws.Range("E14").Value = "Not even rain boots would help at this point."
$ws.Range("F14").Value = "Currently raining very heavily. Raincoat."

$ws.Range("C15").Value = "Raincoat or a Patagonia should do."
$ws.Range("D15").Value = "I'd wear the opposite of a tank top."
$ws.Range("E15").Value = "Are you a polar bear? Then stay inside."
$ws.Range("F15").Value = "Unless you like bad weather, chill indoors."

$ws.Range("C16").Value = "Chilly and a bit wet. You know what to do."
$ws.Range("D16").Value = "Wear a coat, pack an umbrella. You'll need it."
$ws.Range("E16").Value = "Netflix is calling your name. Stay dry."
$ws.Range("F16").Value = "Torrential downpour. No need to leave the house."

# Extra stray phrases off to the side
$ws.Range("J12").Value = "These can be grouped together, probably."
$ws.Range("O12").Value = "It's colder than the polar bear's toenails."

# Merge the note column (F:G) on each phrase row and center it,
# matching the look of the rest of the "Note" column.
$ws.Range("F11:G11").Merge()
$ws.Range("F12:G12").Merge()
$ws.Range("F13:G13").Merge()
$ws.Range("F14:G14").Merge()
$ws.Range("F15:G15").Merge()
$ws.Range("F16:G16").Merge()

$ws.Range("F11:G16").HorizontalAlignment = -4108

# ---------------------------------------------------------------
# Light formatting echoes (border flag carried on these ranges in
# the authored workbook) - apply & release a border so these cells
# pick up a distinct style bucket like the source file.
# ---------------------------------------------------------------
$styleRanges = @("C11","D11","E11","O12","C13","D13","E13","C14","D14","E14","C15","D15","E15","C16","D16","E16","B17","B18","B19","B20")
foreach ($addr in $styleRanges) {
    $ws.Range($addr).Borders.LineStyle = 1
    $ws.Range($addr).Borders.LineStyle = 0
}

$ws.Range("B21").Borders.LineStyle = 1
$ws.Range("B21").Borders.LineStyle = 0
$ws.Range("B22").Borders.LineStyle = 1
$ws.Range("B22").Borders.LineStyle = 0
$ws.Range("B23").Borders.LineStyle = 1
$ws.Range("B23").Borders.LineStyle = 0

# ---------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------
$ws.Range("C:E").ColumnWidth = 51.6640625
$ws.Range("G:G").ColumnWidth = 41.1640625

# ---------------------------------------------------------------
# View state: selection + scrolled position
# ---------------------------------------------------------------
$ws.Range("K17").Select()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
